$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style/format from H1 onto I1:J1 so they match the other headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I and J, rows 2-18
$iVals = @{
    2 = 9
    3 = 8
    4 = 1
    5 = 1
    6 = 1
    7 = 1
    8 = 1
    9 = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 1
    14 = 1
    15 = 1
    16 = 1
    17 = 1
    18 = 1
}

$jVals = @{
    2 = 9
    3 = 8
    4 = 4
    5 = 5
    6 = 5
    7 = 6
    8 = 4
    9 = 4
    10 = 7
    11 = 6
    12 = 5
    13 = 5
    14 = 5
    15 = 3
    16 = 3
    17 = 3
    18 = 2
}

foreach ($r in 2..18) {
    $ws.Cells.Item($r, 9).Value = $iVals[$r]
    $ws.Cells.Item($r, 10).Value = $jVals[$r]
}
